$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44292
$ws.Cells.Item(2, 10).Value = 40
$ws.Cells.Item(2, 11).Value = 3000
$ws.Cells.Item(2, 13).Value = 3000
$ws.Cells.Item(2, 16).Value = 1000

# Row 3
$ws.Cells.Item(3, 4).Value = 44967
$ws.Cells.Item(3, 10).Value = 110
$ws.Cells.Item(3, 11).Value = 3000
$ws.Cells.Item(3, 12).Value = 3300
$ws.Cells.Item(3, 13).Value = 3136
$ws.Cells.Item(3, 16).Value = 1045

# Row 4
$ws.Cells.Item(4, 4).Value = 44389
$ws.Cells.Item(4, 10).Value = 81
$ws.Cells.Item(4, 11).Value = 2800
$ws.Cells.Item(4, 13).Value = 2889
$ws.Cells.Item(4, 16).Value = 963

# Row 5
$ws.Cells.Item(5, 4).Value = 44224
$ws.Cells.Item(5, 10).Value = 67

# Row 6
$ws.Cells.Item(6, 4).Value = 44166
$ws.Cells.Item(6, 11).Value = 2500
$ws.Cells.Item(6, 12).Value = 2500
$ws.Cells.Item(6, 13).Value = 2500
$ws.Cells.Item(6, 16).Value = 833

# Row 7
$ws.Cells.Item(7, 4).Value = 44756
$ws.Cells.Item(7, 10).Value = 104
$ws.Cells.Item(7, 11).Value = 2800
$ws.Cells.Item(7, 12).Value = 3000
$ws.Cells.Item(7, 13).Value = 2904
$ws.Cells.Item(7, 16).Value = 968

# Row 8
$ws.Cells.Item(8, 4).Value = 44557
$ws.Cells.Item(8, 10).Value = 104
$ws.Cells.Item(8, 11).Value = 2000
$ws.Cells.Item(8, 12).Value = 2500
$ws.Cells.Item(8, 13).Value = 2260
$ws.Cells.Item(8, 16).Value = 753

# Row 9
$ws.Cells.Item(9, 4).Value = 44536
$ws.Cells.Item(9, 10).Value = 125
$ws.Cells.Item(9, 11).Value = 2200
$ws.Cells.Item(9, 12).Value = 2200
$ws.Cells.Item(9, 13).Value = 2200
$ws.Cells.Item(9, 16).Value = 733

# Row 10
$ws.Cells.Item(10, 4).Value = 44845
$ws.Cells.Item(10, 12).Value = 2500
$ws.Cells.Item(10, 13).Value = 2500
$ws.Cells.Item(10, 16).Value = 833

# Row 11
$ws.Cells.Item(11, 4).Value = 44260
$ws.Cells.Item(11, 10).Value = 60
$ws.Cells.Item(11, 11).Value = 3500
$ws.Cells.Item(11, 12).Value = 3500
$ws.Cells.Item(11, 13).Value = 3500
$ws.Cells.Item(11, 16).Value = 1167

# Row 12
$ws.Cells.Item(12, 4).Value = 44935
$ws.Cells.Item(12, 10).Value = 78
$ws.Cells.Item(12, 11).Value = 3000
$ws.Cells.Item(12, 13).Value = 3000
$ws.Cells.Item(12, 16).Value = 1000

# Row 13
$ws.Cells.Item(13, 4).Value = 44390
$ws.Cells.Item(13, 10).Value = 50
$ws.Cells.Item(13, 11).Value = 3000
$ws.Cells.Item(13, 12).Value = 3000
$ws.Cells.Item(13, 13).Value = 3000
$ws.Cells.Item(13, 16).Value = 1000

# Row 14
$ws.Cells.Item(14, 4).Value = 44222
$ws.Cells.Item(14, 10).Value = 45

# Row 15
$ws.Cells.Item(15, 4).Value = 44179
$ws.Cells.Item(15, 10).Value = 78
$ws.Cells.Item(15, 11).Value = 3000
$ws.Cells.Item(15, 12).Value = 3000
$ws.Cells.Item(15, 13).Value = 3000
$ws.Cells.Item(15, 16).Value = 1000

# Row 16
$ws.Cells.Item(16, 4).Value = 44225
$ws.Cells.Item(16, 10).Value = 56

# Row 17
$ws.Cells.Item(17, 4).Value = 44937
$ws.Cells.Item(17, 10).Value = 68
$ws.Cells.Item(17, 11).Value = 3500
$ws.Cells.Item(17, 12).Value = 3500
$ws.Cells.Item(17, 13).Value = 3500
$ws.Cells.Item(17, 16).Value = 1167

# Row 18
$ws.Cells.Item(18, 4).Value = 44574
$ws.Cells.Item(18, 10).Value = 50
$ws.Cells.Item(18, 11).Value = 3000
$ws.Cells.Item(18, 12).Value = 3000
$ws.Cells.Item(18, 13).Value = 3000
$ws.Cells.Item(18, 16).Value = 1000

# Row 19
$ws.Cells.Item(19, 4).Value = 44223
$ws.Cells.Item(19, 10).Value = 80
$ws.Cells.Item(19, 11).Value = 2500
$ws.Cells.Item(19, 12).Value = 3000
$ws.Cells.Item(19, 13).Value = 2781
$ws.Cells.Item(19, 16).Value = 927

# Row 20
$ws.Cells.Item(20, 4).Value = 44804
$ws.Cells.Item(20, 10).Value = 85

# Row 21
$ws.Cells.Item(21, 4).Value = 44193
$ws.Cells.Item(21, 10).Value = 70

# Row 22
$ws.Cells.Item(22, 4).Value = 44221
$ws.Cells.Item(22, 10).Value = 50
$ws.Cells.Item(22, 11).Value = 2500
$ws.Cells.Item(22, 13).Value = 2500
$ws.Cells.Item(22, 16).Value = 833

# Row 23
$ws.Cells.Item(23, 4).Value = 44165
$ws.Cells.Item(23, 10).Value = 68

# Row 24
$ws.Cells.Item(24, 4).Value = 44291
$ws.Cells.Item(24, 11).Value = 3000
$ws.Cells.Item(24, 12).Value = 3000
$ws.Cells.Item(24, 13).Value = 3000
$ws.Cells.Item(24, 16).Value = 1000

# Row 25
$ws.Cells.Item(25, 4).Value = 44669
$ws.Cells.Item(25, 10).Value = 92
$ws.Cells.Item(25, 11).Value = 2500
$ws.Cells.Item(25, 12).Value = 3000
$ws.Cells.Item(25, 13).Value = 2755
$ws.Cells.Item(25, 16).Value = 918

# Row 26
$ws.Cells.Item(26, 4).Value = 44242
$ws.Cells.Item(26, 10).Value = 95
$ws.Cells.Item(26, 11).Value = 2500
$ws.Cells.Item(26, 13).Value = 2737
$ws.Cells.Item(26, 16).Value = 912

# Row 27
$ws.Cells.Item(27, 4).Value = 44187
$ws.Cells.Item(27, 10).Value = 65

# Row 28
$ws.Cells.Item(28, 4).Value = 44537
$ws.Cells.Item(28, 10).Value = 88
$ws.Cells.Item(28, 11).Value = 2000
$ws.Cells.Item(28, 12).Value = 2200
$ws.Cells.Item(28, 13).Value = 2091
$ws.Cells.Item(28, 16).Value = 697

# Row 29
$ws.Cells.Item(29, 4).Value = 44243
$ws.Cells.Item(29, 10).Value = 45
$ws.Cells.Item(29, 12).Value = 3000
$ws.Cells.Item(29, 13).Value = 3000
$ws.Cells.Item(29, 16).Value = 1000

# Row 30
$ws.Cells.Item(30, 4).Value = 44340
$ws.Cells.Item(30, 10).Value = 54
$ws.Cells.Item(30, 11).Value = 3000
$ws.Cells.Item(30, 13).Value = 3000
$ws.Cells.Item(30, 16).Value = 1000

# Row 31
$ws.Cells.Item(31, 4).Value = 44627
$ws.Cells.Item(31, 10).Value = 78
$ws.Cells.Item(31, 11).Value = 3500
$ws.Cells.Item(31, 12).Value = 3500
$ws.Cells.Item(31, 13).Value = 3500
$ws.Cells.Item(31, 16).Value = 1167

# Row 32
$ws.Cells.Item(32, 4).Value = 44992
$ws.Cells.Item(32, 10).Value = 45
$ws.Cells.Item(32, 11).Value = 4000
$ws.Cells.Item(32, 12).Value = 4000
$ws.Cells.Item(32, 13).Value = 4000
$ws.Cells.Item(32, 16).Value = 1333

# Row 33
$ws.Cells.Item(33, 4).Value = 44965
$ws.Cells.Item(33, 10).Value = 87
$ws.Cells.Item(33, 11).Value = 3000
$ws.Cells.Item(33, 12).Value = 3000
$ws.Cells.Item(33, 13).Value = 3000
$ws.Cells.Item(33, 16).Value = 1000

# Row 34
$ws.Cells.Item(34, 4).Value = 44559
$ws.Cells.Item(34, 10).Value = 68
$ws.Cells.Item(34, 11).Value = 2000
$ws.Cells.Item(34, 12).Value = 2000
$ws.Cells.Item(34, 13).Value = 2000
$ws.Cells.Item(34, 16).Value = 667
